$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72, pushing the existing row 72 (and anything below) down to 73.
$ws.Rows(72).Insert()

# Fill the newly inserted row 72 with the new data point (2023-07-14 / serial 45121).
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 45121
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = 100112012
$ws.Cells.Item(72, 7).Value = "Espinaca"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 30
$ws.Cells.Item(72, 11).Value = 12000
$ws.Cells.Item(72, 12).Value = 12000
$ws.Cells.Item(72, 13).Value = 12000
$ws.Cells.Item(72, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 1200
$ws.Cells.Item(72, 17).Value = 10
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# Copy the style of column D (date format) from row 73 down into the new row 72's D cell.
$ws.Cells.Item(73, 4).Copy()
$ws.Cells.Item(72, 4).PasteSpecial(-4122)
